$wb = $excel.ActiveWorkbook

# --- Workbook-level: rename the first sheet from "products" to "category" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "category"

# --- Rebuild sheet1 ("category") content -------------------------------
# Remove the old products/paid/quantity/date table entirely, then write the
# new Name/image category table with hyperlinked image URLs.
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "Name"

$names = @("Hand Bags", "Watches", "Shoes", "Suits", "Accessories", "fragrances")
$urls = @(
    "https://images.pexels.com/photos/8801089/pexels-photo-8801089.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1",
    "https://images.pexels.com/photos/190819/pexels-photo-190819.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1",
    "https://images.pexels.com/photos/1598505/pexels-photo-1598505.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1",
    "https://images.pexels.com/photos/1096849/pexels-photo-1096849.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1",
    "https://images.pexels.com/photos/1453008/pexels-photo-1453008.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1",
    "https://images.pexels.com/photos/965989/pexels-photo-965989.jpeg?auto=compress&cs=tinysrgb&w=1260&h=750&dpr=1"
)

# Fill the category names down column A first (matches the original
# authoring order captured in the shared-string table).
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws1.Range("A$row").Value = $names[$i]
}

# Then add the "image" header and the hyperlinked URLs down column B.
$ws1.Range("B1").Value = "image"
for ($i = 0; $i -lt $urls.Count; $i++) {
    $row = $i + 2
    $url = $urls[$i]
    $ws1.Range("B$row").Value = $url
    $ws1.Hyperlinks.Add($ws1.Range("B$row"), $url)
}

# --- Column widths: narrow "Name" column, wide "image" (url) column -----
$ws1.Columns("A").ColumnWidth = 13.88
$ws1.Columns("B").ColumnWidth = 83.02

# --- Page setup: portrait orientation ------------------------------------
$ws1.PageSetup.Orientation = 1

# --- Selection matches the saved view state ------------------------------
$ws1.Range("B12").Select() | Out-Null

Write-Host "category sheet rebuilt"
